$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 2845.6667
$ws.Range("I8").Value = 3380.8
$ws.Range("J8").Value = 170
$ws.Range("K8").Value = 10142.4
$ws.Range("L8").Value = 510
$ws.Range("M8").Value = -10003.4
$ws.Range("N8").Value = -788
$ws.Range("H41").Value = 1200.5294
$ws.Range("I41").Value = 1161.4
$ws.Range("J41").Value = 1256.4286
$ws.Range("K41").Value = 1161.4
$ws.Range("L41").Value = 1256.4286
$ws.Range("M41").Value = -721.4000000000001
$ws.Range("N41").Value = -2136.4286
$ws.Range("H62").Value = 4003.2
$ws.Range("I62").Value = 4003.3333
$ws.Range("J62").Value = 4003
$ws.Range("K62").Value = 4003.3333
$ws.Range("L62").Value = 4003
$ws.Range("M62").Value = -3379.3333
$ws.Range("N62").Value = -5251
$ws.Range("H65").Value = 4003.2
$ws.Range("I65").Value = 4003.3333
$ws.Range("J65").Value = 4003
$ws.Range("K65").Value = 20016.6665
$ws.Range("L65").Value = 20015
$ws.Range("M65").Value = -16896.6665
$ws.Range("N65").Value = -26255
$ws.Range("H98").Value = 100029460
$ws.Range("I98").Value = 100029460
$ws.Range("K98").Value = 100029460
$ws.Range("M98").Value = -100027962
$ws.Range("H112").Value = 4750
$ws.Range("J112").Value = 6000
$ws.Range("L112").Value = 18000
$ws.Range("N112").Value = -20216
$ws.Range("H122").Value = 100029460
$ws.Range("I122").Value = 100029460
$ws.Range("K122").Value = 300088380
$ws.Range("M122").Value = -300085930
$ws.Range("H135").Value = 5379.143
$ws.Range("I135").Value = 648.1
$ws.Range("K135").Value = 5832.900000000001
$ws.Range("M135").Value = -3297.900000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1376.6471
$ws.Range("I2").Value = 1446.9333
$ws.Range("K2").Value = 1446.9333
$ws.Range("M2").Value = -1333.9333
$ws.Range("H32").Value = 8475759
$ws.Range("I32").Value = 8475759
$ws.Range("K32").Value = 8475759
$ws.Range("M32").Value = -8475472
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H45").Value = 2393.4583
$ws.Range("J45").Value = 2346.1333
$ws.Range("L45").Value = 2346.1333
$ws.Range("N45").Value = -3100.1333
$ws.Range("H96").Value = 28003.1
$ws.Range("J96").Value = 28003.1
$ws.Range("L96").Value = 28003.1
$ws.Range("N96").Value = -33495.1
$ws.Range("H97").Value = 952.25
$ws.Range("I97").Value = 982
$ws.Range("K97").Value = 982
$ws.Range("M97").Value = -486
$ws.Range("H116").Value = 1376.6471
$ws.Range("I116").Value = 1446.9333
$ws.Range("K116").Value = 1446.9333
$ws.Range("M116").Value = 847.0667000000001
$ws.Range("H122").Value = 1238.5714
$ws.Range("I122").Value = 930
$ws.Range("K122").Value = 2790
$ws.Range("M122").Value = -340
$ws.Range("H132").Value = 11085
$ws.Range("I132").Value = 6429.875
$ws.Range("K132").Value = 19289.625
$ws.Range("M132").Value = -16759.625
$ws.Range("J139").Value = 59000
$ws.Range("L139").Value = 59000
$ws.Range("N139").Value = -69280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1376.6471
$ws.Range("I3").Value = 1446.9333
$ws.Range("K3").Value = 1446.9333
$ws.Range("M3").Value = -1332.9333
$ws.Range("H20").Value = 4731.143
$ws.Range("I20").Value = 5286.3335
$ws.Range("K20").Value = 5286.3335
$ws.Range("M20").Value = -5039.3335
$ws.Range("H38").Value = 39001
$ws.Range("J38").Value = 39001
$ws.Range("L38").Value = 39001
$ws.Range("N38").Value = -39833
$ws.Range("H134").Value = 68494
$ws.Range("I134").Value = 3344.1428
$ws.Range("J134").Value = 119166.11
$ws.Range("K134").Value = 10032.4284
$ws.Range("L134").Value = 357498.33
$ws.Range("M134").Value = -7497.428400000001
$ws.Range("N134").Value = -362568.33

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 939.73334
$ws.Range("I22").Value = 989
$ws.Range("K22").Value = 989
$ws.Range("M22").Value = -639
$ws.Range("H31").Value = 1236693.2
$ws.Range("I31").Value = 16692.1
$ws.Range("K31").Value = 16692.1
$ws.Range("M31").Value = -16397.1
$ws.Range("H34").Value = 1236693.2
$ws.Range("I34").Value = 16692.1
$ws.Range("K34").Value = 16692.1
$ws.Range("M34").Value = -16490.1
$ws.Range("H58").Value = 3778.1177
$ws.Range("I58").Value = 3670.6924
$ws.Range("K58").Value = 3670.6924
$ws.Range("M58").Value = -3467.6924
$ws.Range("H99").Value = 2095.35
$ws.Range("I99").Value = 1979.8823
$ws.Range("K99").Value = 1979.8823
$ws.Range("M99").Value = -481.8823
$ws.Range("H126").Value = 2095.35
$ws.Range("I126").Value = 1979.8823
$ws.Range("K126").Value = 5939.6469
$ws.Range("M126").Value = -3469.6469
$ws.Range("H136").Value = 3778.1177
$ws.Range("I136").Value = 3670.6924
$ws.Range("K136").Value = 11012.0772
$ws.Range("M136").Value = -8462.0772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 82997.5
$ws.Range("J37").Value = 82997.5
$ws.Range("L37").Value = 248992.5
$ws.Range("N37").Value = -249216.5
$ws.Range("H119").Value = 5935.75
$ws.Range("I119").Value = 1403.625
$ws.Range("K119").Value = 4210.875
$ws.Range("M119").Value = 627.125
$ws.Range("H122").Value = 1350.3334
$ws.Range("I122").Value = 866.6667
$ws.Range("J122").Value = 1592.1666
$ws.Range("K122").Value = 7800.0003
$ws.Range("L122").Value = 14329.4994
$ws.Range("M122").Value = -5350.0003
$ws.Range("N122").Value = -19229.4994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4112.5
$ws.Range("I102").Value = 3472.8696
$ws.Range("K102").Value = 3472.8696
$ws.Range("M102").Value = -1850.8696
$ws.Range("H141").Value = 20000
$ws.Range("J141").Value = 20000
$ws.Range("L141").Value = 20000
$ws.Range("N141").Value = -30360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 45454828
$ws.Range("I55").Value = 58823844
$ws.Range("K55").Value = 58823844
$ws.Range("M55").Value = -58823671
$ws.Range("H101").Value = 65227
$ws.Range("J101").Value = 65227
$ws.Range("L101").Value = 65227
$ws.Range("N101").Value = -71717
$ws.Range("H122").Value = 4538.8887
$ws.Range("I122").Value = 3680
$ws.Range("J122").Value = 8833.333000000001
$ws.Range("K122").Value = 11040
$ws.Range("L122").Value = 26499.999
$ws.Range("M122").Value = -8590
$ws.Range("N122").Value = -31399.999
$ws.Range("H132").Value = 589017.25
$ws.Range("I132").Value = 1002142.3
$ws.Range("J132").Value = 129989.445
$ws.Range("K132").Value = 3006426.9
$ws.Range("L132").Value = 389968.335
$ws.Range("M132").Value = -3003896.9
$ws.Range("N132").Value = -395028.335
$ws.Range("H136").Value = 114019.3
$ws.Range("I136").Value = 3773.5
$ws.Range("K136").Value = 11320.5
$ws.Range("M136").Value = -8770.5
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 33338912
$ws.Range("J62").Value = 50004424
$ws.Range("L62").Value = 50004424
$ws.Range("N62").Value = -50005672
$ws.Range("H65").Value = 33338912
$ws.Range("J65").Value = 50004424
$ws.Range("L65").Value = 250022120
$ws.Range("N65").Value = -250028360
$ws.Range("H81").Value = 1491.75
$ws.Range("I81").Value = 1491.75
$ws.Range("K81").Value = 2983.5
$ws.Range("M81").Value = -1922.5
$ws.Range("H84").Value = 1491.75
$ws.Range("I84").Value = 1491.75
$ws.Range("K84").Value = 14917.5
$ws.Range("M84").Value = -9613.5
$ws.Range("H96").Value = 3888.4443
$ws.Range("J96").Value = 7666.3335
$ws.Range("L96").Value = 7666.3335
$ws.Range("N96").Value = -10412.3335
$ws.Range("H98").Value = 44594
$ws.Range("J98").Value = 44594
$ws.Range("L98").Value = 44594
$ws.Range("N98").Value = -50584
$ws.Range("H136").Value = 15832.823
$ws.Range("I136").Value = 1596.5
$ws.Range("J136").Value = 50000
$ws.Range("K136").Value = 4789.5
$ws.Range("L136").Value = 150000
$ws.Range("M136").Value = -2239.5
$ws.Range("N136").Value = -155100
$ws.Range("H140").Value = 60000
$ws.Range("J140").Value = 60000
$ws.Range("L140").Value = 60000
$ws.Range("N140").Value = -70360
